# New crime data collected
# Updates the CompStat_1 sheet: header "Volume/Number" + report week dates,
# and the weekly crime-complaint figures in rows 14-27 (plus the two N/A
# placeholder cells in rows 28-29 that revert to blank/"***.*" markers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: set a plain numeric value onto a cell that is ALREADY numeric
# (keeps its existing number format / style untouched).
# ---------------------------------------------------------------------
function Set-NumVal($addr, $val) {
  $ws.Range($addr).Value2 = $val
}

# ---------------------------------------------------------------------
# Helper: convert a cell that currently holds the text placeholder
# ("N/A" / "***.*") into a real number, adopting the number format of a
# sibling cell ($refAddr) that already carries the correct numeric style.
# ---------------------------------------------------------------------
function Set-NumCell($addr, $val, $refAddr) {
  $c = $ws.Range($addr)
  $ref = $ws.Range($refAddr)
  $c.NumberFormat = $ref.NumberFormat
  $c.Value2 = $val
}

# ---------------------------------------------------------------------
# Helper: convert a cell that currently holds a number back into the
# text placeholder, reusing the formatting (general number format /
# right-aligned text style) of a sibling cell ($refAddr) that already
# has that placeholder style.
# ---------------------------------------------------------------------
function Set-TextCell($addr, $text, $refAddr) {
  $c = $ws.Range($addr)
  $ref = $ws.Range($refAddr)
  $c.NumberFormat = "@"
  $c.Value2 = $text
  $ref.Copy()
  $c.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# ---------------------------------------------------------------------
# Header: "Volume 31   Number  7" -> "...  8"  (A8, rich text run #4)
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 1).Text = "8"

# ---------------------------------------------------------------------
# Header: report week "2/12/2024 .. 2/18/2024" -> "2/19/2024 .. 2/25/2024"
# (C9, rich text runs #2 and #4)
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(27, 9).Text = "2/19/2024"
$ws.Range("C9").Characters(47, 9).Text = "2/25/2024"

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
Set-NumCell "D14" 1 "I14"
Set-NumCell "E14" -100 "K14"
Set-NumCell "G14" 1 "I14"
Set-NumCell "H14" -100 "K14"
Set-NumVal "J14" 2
Set-NumVal "K14" -50
Set-NumVal "N14" -85.714285714285

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-TextCell "C15" "0" "D15"
Set-NumVal "F15" 2
Set-NumVal "H15" 100

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-NumVal "C16" 11
Set-NumVal "E16" 10
Set-NumVal "F16" 44
Set-NumVal "G16" 32
Set-NumVal "H16" 37.5
Set-NumVal "I16" 76
Set-NumVal "J16" 56
Set-NumVal "K16" 35.714285714285
Set-NumVal "L16" 65.217391304347
Set-NumVal "M16" 13.432835820895
Set-NumVal "N16" -71.747211895910

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-NumVal "C17" 8
Set-NumVal "D17" 18
Set-NumVal "E17" -55.555555555555
Set-NumVal "F17" 41
Set-NumVal "G17" 50
Set-NumVal "H17" -18
Set-NumVal "I17" 87
Set-NumVal "J17" 100
Set-NumVal "K17" -13
Set-NumVal "L17" 8.75
Set-NumVal "M17" 163.636363636364
Set-NumVal "N17" 52.631578947368

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Set-NumVal "C18" 2
Set-NumVal "D18" 3
Set-NumVal "E18" -33.333333333333
Set-NumVal "F18" 23
Set-NumVal "G18" 20
Set-NumVal "H18" 15
Set-NumVal "I18" 51
Set-NumVal "J18" 38
Set-NumVal "K18" 34.210526315789
Set-NumVal "L18" 131.818181818182
Set-NumVal "M18" -25
Set-NumVal "N18" -87.560975609756

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
Set-NumVal "C19" 26
Set-NumVal "D19" 24
Set-NumVal "E19" 8.333333333333
Set-NumVal "F19" 86
Set-NumVal "G19" 63
Set-NumVal "H19" 36.507936507936
Set-NumVal "I19" 171
Set-NumVal "J19" 121
Set-NumVal "K19" 41.322314049586
Set-NumVal "L19" -38.043478260869
Set-NumVal "M19" 81.914893617021
Set-NumVal "N19" 0

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-NumVal "C20" 3
Set-NumVal "D20" 5
Set-NumVal "E20" -40
Set-NumVal "F20" 16
Set-NumVal "G20" 20
Set-NumVal "H20" -20
Set-NumVal "I20" 40
Set-NumVal "J20" 42
Set-NumVal "K20" -4.761904761904
Set-NumVal "L20" 33.333333333333
Set-NumVal "M20" 73.913043478260
Set-NumVal "N20" -88.826815642458

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
Set-NumVal "C21" 50
Set-NumVal "D21" 61
Set-NumVal "E21" -18.032786885245
Set-NumVal "F21" 212
Set-NumVal "G21" 187
Set-NumVal "H21" 13.368983957219
Set-NumVal "I21" 432
Set-NumVal "J21" 362
Set-NumVal "K21" 19.337016574585
Set-NumVal "L21" -5.882352941176
Set-NumVal "M21" 48.965517241379
Set-NumVal "N21" -66.064414768263

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-TextCell "C22" "0" "C23"
Set-NumCell "D22" 2 "F22"
Set-NumCell "E22" -100 "H22"
Set-NumVal "F22" 2
Set-NumVal "G22" 4
Set-NumVal "H22" -50
Set-NumVal "J22" 5
Set-NumVal "K22" 80

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-NumVal "C24" 74
Set-NumVal "D24" 70
Set-NumVal "E24" 5.714285714285
Set-NumVal "F24" 267
Set-NumVal "G24" 240
Set-NumVal "H24" 11.25
Set-NumVal "I24" 522
Set-NumVal "J24" 480
Set-NumVal "K24" 8.75
Set-NumVal "L24" 43.406593406593
Set-NumVal "M24" 141.666666666667

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
Set-NumVal "C25" 29
Set-NumVal "E25" 52.631578947368
Set-NumVal "F25" 101
Set-NumVal "G25" 66
Set-NumVal "H25" 53.030303030303
Set-NumVal "I25" 184
Set-NumVal "J25" 141
Set-NumVal "K25" 30.496453900709
Set-NumVal "L25" 87.755102040816
Set-NumVal "M25" 127.16049382716

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
Set-TextCell "C26" "0" "D15"
Set-NumVal "E26" -100
Set-NumVal "F26" 2
Set-NumVal "G26" 4
Set-NumVal "H26" -50
Set-NumVal "J26" 6
Set-NumVal "K26" 0

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
Set-NumVal "C27" 2
Set-NumVal "D27" 5
Set-NumVal "E27" -60
Set-NumVal "F27" 13
Set-NumVal "G27" 10
Set-NumVal "H27" 30
Set-NumVal "I27" 21
Set-NumVal "J27" 23
Set-NumVal "K27" -8.695652173913
Set-NumVal "L27" 61.538461538461

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic. (goes back to "no activity this week")
# ---------------------------------------------------------------------
Set-TextCell "D28" "0" "C28"
Set-TextCell "E28" "***.*" "C23"

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc. (goes back to "no activity this week")
# ---------------------------------------------------------------------
Set-TextCell "D29" "0" "C28"
Set-TextCell "E29" "***.*" "C23"
